# Bump the "Förändrad" (changed) date column (C) for every data row
# from 2023-09-11 (45180) to 2023-09-12 (45181).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C408").Value = 45181
